# Publication release 0.2.0 preparation
# - bump Version property
# - bump Date property
# - add a new "Jurisdiction" property row (iso:code:3166:FR) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (A3/B3) and Date (A8/B8) values in place.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new row right after "Contact" (row 10), before "Description" (row 11),
# shifting everything below it down by one row, and copy the existing row
# formatting down onto it so the new row keeps the same borders / wrap text
# style as the rest of the table.
$ws.Range("A11:B11").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
